# Connected Office Test Data - "Added Functionality to save test results"
#
# The "Test Results" worksheet tracks CRUD test pass/fail flags per Zone.
# The new "save test results" functionality flips the Create Test Passed
# flag (column B) for zones Z06-Z09 (rows 7-10) from TRUE to FALSE, and
# leaves the current selection on G9 (where the user last clicked while
# exercising the save feature).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Create Test Passed" results for Z06 (row7) .. Z09 (row10)
$ws.Range("B7").Value = $false
$ws.Range("B8").Value = $false
$ws.Range("B9").Value = $false
$ws.Range("B10").Value = $false

# Move the active selection to reflect where the user left off
$ws.Range("G9").Select()
